$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns with latest values.
# D-column values are numeric-looking strings (e.g. "1.00", "7.42") that must stay
# text, matching the workbook convention (thousands-dot formatted prices). Force
# text interpretation via NumberFormat "@", then restore the original General style
# (copied from the untouched "B" cell in the same row) so no stray style id lingers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.047.56'
$ws.Range('D2').Style = $ws.Range('B2').Style
$ws.Range('E2').Value = '  -3.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.284.80'
$ws.Range('D3').Style = $ws.Range('B3').Style
$ws.Range('E3').Value = '  -4.05%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '554.96'
$ws.Range('D5').Style = $ws.Range('B5').Style
$ws.Range('E5').Value = '  -4.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.74'
$ws.Range('D6').Style = $ws.Range('B6').Style
$ws.Range('E6').Value = '  -8.31%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.285.94'
$ws.Range('D8').Style = $ws.Range('B8').Style
$ws.Range('E8').Value = '  -4.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.465'
$ws.Range('D9').Style = $ws.Range('B9').Style
$ws.Range('E9').Value = '  -3.95%  '
$ws.Range('E10').Value = '  -2.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.118'
$ws.Range('D11').Style = $ws.Range('B11').Style
$ws.Range('E11').Value = '  -5.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.405'
$ws.Range('D12').Style = $ws.Range('B12').Style
$ws.Range('E12').Value = '  -3.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.842.58'
$ws.Range('D13').Style = $ws.Range('B13').Style
$ws.Range('E13').Value = '  -4.19%  '
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.65'
$ws.Range('D15').Style = $ws.Range('B15').Style
$ws.Range('E15').Value = '  -6.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.277.91'
$ws.Range('D16').Style = $ws.Range('B16').Style
$ws.Range('E16').Value = '  -4.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000163'
$ws.Range('D17').Style = $ws.Range('B17').Style
$ws.Range('E17').Value = '  -5.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '60.094.48'
$ws.Range('D18').Style = $ws.Range('B18').Style
$ws.Range('E18').Value = '  -3.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.04'
$ws.Range('D19').Style = $ws.Range('B19').Style
$ws.Range('E19').Value = '  -7.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.73'
$ws.Range('D20').Style = $ws.Range('B20').Style
$ws.Range('E20').Value = '  -5.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.48'
$ws.Range('D21').Style = $ws.Range('B21').Style
$ws.Range('E21').Value = '  -5.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '371.68'
$ws.Range('D22').Style = $ws.Range('B22').Style
$ws.Range('E22').Value = '  -2.88%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.10'
$ws.Range('D24').Style = $ws.Range('B24').Style
$ws.Range('E24').Value = '  -4.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.531'
$ws.Range('D25').Style = $ws.Range('B25').Style
$ws.Range('E25').Value = '  -7.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.408.72'
$ws.Range('D26').Style = $ws.Range('B26').Style
$ws.Range('E26').Value = '  -4.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000103'
$ws.Range('D27').Style = $ws.Range('B27').Style
$ws.Range('E27').Value = '  -8.81%  '
$ws.Range('E28').Value = '  -3.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = $ws.Range('B29').Style
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.05'
$ws.Range('D30').Style = $ws.Range('B30').Style
$ws.Range('E30').Value = '  -7.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = $ws.Range('B31').Style
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  -5.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.42'
$ws.Range('D33').Style = $ws.Range('B33').Style
$ws.Range('E33').Value = '  -5.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '22.50'
$ws.Range('D34').Style = $ws.Range('B34').Style
$ws.Range('E34').Value = '  -3.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.23'
$ws.Range('D35').Style = $ws.Range('B35').Style
$ws.Range('E35').Value = '  -7.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '166.14'
$ws.Range('D36').Style = $ws.Range('B36').Style
$ws.Range('E36').Value = '  -1.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.04'
$ws.Range('D37').Style = $ws.Range('B37').Style
$ws.Range('E37').Value = '  -9.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.52'
$ws.Range('D38').Style = $ws.Range('B38').Style
$ws.Range('E38').Value = '  -5.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.60'
$ws.Range('D39').Style = $ws.Range('B39').Style
$ws.Range('E39').Value = '  -5.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.312.46'
$ws.Range('D40').Style = $ws.Range('B40').Style
$ws.Range('E40').Value = '  -4.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0721'
$ws.Range('D41').Style = $ws.Range('B41').Style
$ws.Range('E41').Value = '  -8.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '25.68'
$ws.Range('D42').Style = $ws.Range('B42').Style
$ws.Range('E42').Value = '  -17.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.52'
$ws.Range('D43').Style = $ws.Range('B43').Style
$ws.Range('E43').Value = '  -2.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.744'
$ws.Range('D44').Style = $ws.Range('B44').Style
$ws.Range('E44').Value = '  -4.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.09'
$ws.Range('D45').Style = $ws.Range('B45').Style
$ws.Range('E45').Value = '  -8.04%  '
$ws.Range('E46').Value = '  -4.15%  '
$ws.Range('E47').Value = '  -7.16%  '
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.311.41'
$ws.Range('D49').Style = $ws.Range('B49').Style
$ws.Range('E49').Value = '  -9.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.31'
$ws.Range('D50').Style = $ws.Range('B50').Style
$ws.Range('E50').Value = '  -7.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.37'
$ws.Range('D51').Style = $ws.Range('B51').Style
$ws.Range('E51').Value = '  -5.69%  '
